$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row restyle (swap which header cells carry which bold style) ---
$ws.Range("A1").Style = $ws.Range("C1").Style
$ws.Range("A2").Style = $ws.Range("C2").Style

# --- New font applied to the "No response" row labels ---
$ws.Range("A6").Font.Name = "Times Roman"
$ws.Range("A7").Font.Name = "Times Roman"

# --- Study 2 demographic counts replacing Study 1 percentages ---
$ws.Range("B14").Value = "196 (68.06%)"
$ws.Range("C14").Value = "30 (65.22%)"
$ws.Range("B15").Value = "27 (9.38)%"
$ws.Range("C15").Value = "4 (8.7%)"
$ws.Range("B16").Value = "0 (0%)"
$ws.Range("C16").Value = "0 (0%)"
$ws.Range("B17").Value = "22 (7.64%)"
$ws.Range("C17").Value = "4 (8.7%)"
$ws.Range("B18").Value = "0 (0%)"
$ws.Range("C18").Value = "0 (0%)"
$ws.Range("B19").Value = "16 (5.56%)"
$ws.Range("C19").Value = "3 (6.52%)"
$ws.Range("B20").Value = "25 (8.68%)"
$ws.Range("C20").Value = "4 (8.7%)"
$ws.Range("B21").Value = "2 (0.69%)"
$ws.Range("C21").Value = "1 (2.17%)"

# --- Selection position, matching the author's cursor when saving ---
$ws.Range("E20").Select()
